# Apply 2024-03-27 crime data update across Citywide Totals, By Neighborhood, and all
# per-neighborhood sheets. Each block targets one worksheet by its tab position (1-based),
# which matches the workbook's physical sheetN.xml ordering, to avoid any ambiguity with
# sheet names that contain special characters (e.g. "Rush & Division").

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value2 = 1665
$ws.Range("K3").Value2 = 1577
$ws.Range("E4").Value2 = 2028
$ws.Range("I4").Value2 = 1785
$ws.Range("J4").Value2 = 1796
$ws.Range("K4").Value2 = 342
$ws.Range("K5").Value2 = 100
$ws.Range("K6").Value2 = 2053
$ws.Range("E7").Value2 = 26033
$ws.Range("I7").Value2 = 26239
$ws.Range("J7").Value2 = 29265
$ws.Range("K7").Value2 = 5737

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("K7").Value2 = 161
$ws.Range("K8").Value2 = 368
$ws.Range("K10").Value2 = 34
$ws.Range("K11").Value2 = 120
$ws.Range("K15").Value2 = 50
$ws.Range("K18").Value2 = 44
$ws.Range("K19").Value2 = 156
$ws.Range("K20").Value2 = 130
$ws.Range("K25").Value2 = 28
$ws.Range("K29").Value2 = 260
$ws.Range("K30").Value2 = 22
$ws.Range("K33").Value2 = 236
$ws.Range("K37").Value2 = 193
$ws.Range("K42").Value2 = 205
$ws.Range("K44").Value2 = 52
$ws.Range("K46").Value2 = 11
$ws.Range("K48").Value2 = 60
$ws.Range("K52").Value2 = 151
$ws.Range("K54").Value2 = 97
$ws.Range("K55").Value2 = 63
$ws.Range("J63").Value2 = 91
$ws.Range("K63").Value2 = 20
$ws.Range("K65").Value2 = 145
$ws.Range("K66").Value2 = 24
$ws.Range("K67").Value2 = 224
$ws.Range("K79").Value2 = 154
$ws.Range("K80").Value2 = 19
$ws.Range("K83").Value2 = 118
$ws.Range("E84").Value2 = 170
$ws.Range("I84").Value2 = 226
$ws.Range("K84").Value2 = 40
$ws.Range("K85").Value2 = 291
$ws.Range("K86").Value2 = 39
$ws.Range("K88").Value2 = 70
$ws.Range("K89").Value2 = 75
$ws.Range("K91").Value2 = 52
$ws.Range("K93").Value2 = 25
$ws.Range("K94").Value2 = 67
$ws.Range("K95").Value2 = 99
$ws.Range("K96").Value2 = 78
$ws.Range("K97").Value2 = 49
$ws.Range("E101").Value2 = 26033
$ws.Range("I101").Value2 = 26239
$ws.Range("J101").Value2 = 29265
$ws.Range("K101").Value2 = 5737

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item(3)
$ws.Range("K2").Value2 = 14
$ws.Range("K6").Value2 = 12

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item(4)
$ws.Range("K2").Value2 = 28
$ws.Range("K7").Value2 = 78

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("K2").Value2 = 60
$ws.Range("K3").Value2 = 52
$ws.Range("K7").Value2 = 161

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("K6").Value2 = 49
$ws.Range("K7").Value2 = 120

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Range("K3").Value2 = 26
$ws.Range("K7").Value2 = 75

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K2").Value2 = 106
$ws.Range("K3").Value2 = 94
$ws.Range("K7").Value2 = 291

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item(9)
$ws.Range("K3").Value2 = 32
$ws.Range("K6").Value2 = 74
$ws.Range("K7").Value2 = 151

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K2").Value2 = 110
$ws.Range("K3").Value2 = 106
$ws.Range("K5").Value2 = 8
$ws.Range("K6").Value2 = 125
$ws.Range("K7").Value2 = 368

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K2").Value2 = 50
$ws.Range("K3").Value2 = 39
$ws.Range("K6").Value2 = 22
$ws.Range("K7").Value2 = 118

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Range("K3").Value2 = 88
$ws.Range("K4").Value2 = 15
$ws.Range("K7").Value2 = 236

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Range("K2").Value2 = 35
$ws.Range("K7").Value2 = 99

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K2").Value2 = 42
$ws.Range("K7").Value2 = 193

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K3").Value2 = 35
$ws.Range("K6").Value2 = 64
$ws.Range("K7").Value2 = 145

# Sheet 19: Fuller Park
$ws = $wb.Worksheets.Item(19)
$ws.Range("K2").Value2 = 6
$ws.Range("K7").Value2 = 22

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Range("K4").Value2 = 11
$ws.Range("K7").Value2 = 224

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item(22)
$ws.Range("K2").Value2 = 15
$ws.Range("K3").Value2 = 12
$ws.Range("E4").Value2 = 15
$ws.Range("I4").Value2 = 6
$ws.Range("K6").Value2 = 11
$ws.Range("E7").Value2 = 170
$ws.Range("I7").Value2 = 226
$ws.Range("K7").Value2 = 40

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("K2").Value2 = 19
$ws.Range("K6").Value2 = 38
$ws.Range("K7").Value2 = 97

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value2 = 70
$ws.Range("K3").Value2 = 84
$ws.Range("K7").Value2 = 260

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K2").Value2 = 14
$ws.Range("K7").Value2 = 60

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K3").Value2 = 48
$ws.Range("K6").Value2 = 50
$ws.Range("K7").Value2 = 156

# Sheet 28: Irving Park
$ws = $wb.Worksheets.Item(28)
$ws.Range("K2").Value2 = 8
$ws.Range("K7").Value2 = 52

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("K6").Value2 = 91
$ws.Range("K7").Value2 = 205

# Sheet 34: Avondale
$ws = $wb.Worksheets.Item(34)
$ws.Range("K3").Value2 = 4
$ws.Range("K7").Value2 = 34

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("K2").Value2 = 26
$ws.Range("K7").Value2 = 63

# Sheet 38: Jefferson Park
$ws = $wb.Worksheets.Item(38)
$ws.Range("K3").Value2 = 3
$ws.Range("K7").Value2 = 11

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Range("K2").Value2 = 17
$ws.Range("K7").Value2 = 52

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("K3").Value2 = 54
$ws.Range("K7").Value2 = 154

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("K2").Value2 = 37
$ws.Range("K5").Value2 = 1
$ws.Range("K7").Value2 = 130

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Range("K3").Value2 = 13
$ws.Range("K7").Value2 = 44

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item(48)
$ws.Range("K6").Value2 = 13
$ws.Range("K7").Value2 = 25

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K6").Value2 = 29
$ws.Range("K7").Value2 = 67

# Sheet 52: East Side
$ws = $wb.Worksheets.Item(52)
$ws.Range("K2").Value2 = 11
$ws.Range("K3").Value2 = 11
$ws.Range("K7").Value2 = 28

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item(54)
$ws.Range("K3").Value2 = 10
$ws.Range("K7").Value2 = 50

# Sheet 59: North Center
$ws = $wb.Worksheets.Item(59)
$ws.Range("K6").Value2 = 11
$ws.Range("K7").Value2 = 24

# Sheet 65: West Town
$ws = $wb.Worksheets.Item(65)
$ws.Range("K6").Value2 = 32
$ws.Range("K7").Value2 = 49

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Range("K6").Value2 = 40
$ws.Range("K7").Value2 = 70

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item(72)
$ws.Range("K3").Value2 = 8
$ws.Range("K7").Value2 = 39

# Sheet 87: Rush &amp; Division
$ws = $wb.Worksheets.Item(87)
$ws.Range("K6").Value2 = 9
$ws.Range("K7").Value2 = 19
